# Fix: use MP (monetary-policy) surprises as the correct shock.
# The "home_ownership_pct" series in column J had been built one quarter
# ahead of where it should line up; shift the whole J6:J65 block down by
# one year (4 quarterly rows) so each block of 4 quarters gets the value
# that previously belonged to the preceding block. Rows J2:J5 (the first
# year) stay as-is, and the newly-revealed tail (J62:J65) is filled with
# the value that used to sit in J58:J61.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for J2:J61 first, since we will overwrite
# J6:J65 using data read from 4 rows above - reading directly off the
# sheet while writing could clobber source data before it's copied.
$sourceValues = @{}
for ($r = 2; $r -le 61; $r++) {
    $sourceValues[$r] = $ws.Cells.Item($r, 10).Value2
}

for ($r = 65; $r -ge 6; $r--) {
    $ws.Cells.Item($r, 10).Value = $sourceValues[$r - 4]
}
